$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.985.51"
$ws.Range("E2").Value = "  -3.82%  "
$ws.Range("D3").Value = "2.527.85"
$ws.Range("E3").Value = "  -4.48%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'545.20"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'147.39"
$ws.Range("E6").Value = "  -4.56%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "2.525.67"
$ws.Range("E9").Value = "  -4.54%  "
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").Value = "'5.51"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").Value = "2.938.74"
$ws.Range("E14").Value = "  -5.66%  "
$ws.Range("D15").Value = "'24.63"
$ws.Range("E15").Value = "  -4.30%  "
$ws.Range("D16").Value = "59.926.61"
$ws.Range("E16").Value = "  -3.75%  "
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "2.508.94"
$ws.Range("E18").Value = "  -5.45%  "
$ws.Range("D19").Value = "'11.42"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "'4.38"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").Value = "'327.40"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("D22").Value = "'0.993"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").Value = "'61.42"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").Value = "'0.451"
$ws.Range("E25").Value = "  -10.30%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").Value = "'7.83"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "0.0₃0803"
$ws.Range("D30").Value = "'1.30"
$ws.Range("E30").Value = "  -5.09%  "
$ws.Range("D31").Value = "'6.81"
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").Value = "'0.996"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'157.84"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'1.45"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").Value = "'18.99"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "'4.54"
$ws.Range("E37").Value = "  -4.40%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").Value = "'5.89"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("D40").Value = "'317.89"
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.80"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'36.85"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("D43").Value = "'0.841"
$ws.Range("E43").Value = "  -7.44%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'0.605"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").Value = "'10.70"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").Value = "'126.90"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Value = "'0.0536"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").Value = "'0.0943"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'18.76"
$ws.Range("E51").Value = "  -5.38%  "

# Reset style on touched cells so forced-text (quote-prefix) entries do not
# pick up a NumberFormat/quotePrefix style that the source file never had.
$ws.Range("D2:E51").Style = "Normal"
$ws.Range("B41:C42").Style = "Normal"
$ws.Range("B51:C51").Style = "Normal"
